# Insert a new data row at row 189 ("Fruta / hortaliza, semanal" update),
# pushing the existing rows 189-253 down to 190-254.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("189:189").Insert()

$ws.Range("A189").Value = 1
$ws.Range("B189").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C189").Value = "Arica y Parinacota"
$ws.Range("D189").Value = 44722
$ws.Range("E189").Value = 15
$ws.Range("F189").Value = "Fruta"
$ws.Range("G189").Value = 100108
$ws.Range("H189").Value = "Tropicales y subtropicales"
$ws.Range("I189").Value = 100108006
$ws.Range("J189").Value = "Plátano"
$ws.Range("K189").Value = "Sin especificar"
$ws.Range("L189").Value = "Pintón"
$ws.Range("M189").Value = 120
$ws.Range("N189").Value = 13000
$ws.Range("O189").Value = 14000
$ws.Range("P189").Value = 13500
$ws.Range("Q189").Value = "`$/caja 20 kilos"
$ws.Range("R189").Value = "Ecuador"
$ws.Range("S189").Value = 675
$ws.Range("T189").Value = 20
